$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-parsed as a number by
# Excel (e.g. "1.010" -> 1.01, "54.00" -> 54). Force them to Text format
# first so the literal string is preserved exactly, matching the source data.
$forceTextCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.062.52'
$ws.Range("E2").Value = '  -4.04%  '

$ws.Range("D3").Value = '1.960.21'
$ws.Range("E3").Value = '  -6.57%  '

$ws.Range("D4").Value = '1.010'
$ws.Range("E4").Value = '  +0.80%  '

$ws.Range("D5").Value = '327.09'
$ws.Range("E5").Value = '  -4.54%  '

$ws.Range("D6").Value = '1.009'
$ws.Range("E6").Value = '  +0.76%  '

$ws.Range("D7").Value = '0.4972'
$ws.Range("E7").Value = '  -5.86%  '

$ws.Range("D8").Value = '0.4201'
$ws.Range("E8").Value = '  -4.30%  '

$ws.Range("D9").Value = '54.00'
$ws.Range("E9").Value = '  -1.97%  '

$ws.Range("D10").Value = '0.08978'
$ws.Range("E10").Value = '  -4.32%  '

$ws.Range("D11").Value = '1.096'
$ws.Range("E11").Value = '  -6.72%  '

$ws.Range("D12").Value = '22.90'
$ws.Range("E12").Value = '  -7.54%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '7.846'
$ws.Range("E13").Value = '  -8.39%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.940.60'
$ws.Range("E14").Value = '  -6.34%  '

$ws.Range("D15").Value = '6.413'
$ws.Range("E15").Value = '  -6.61%  '

$ws.Range("D16").Value = '1.011'
$ws.Range("E16").Value = '  +0.91%  '

$ws.Range("D17").Value = '0.00001094'
$ws.Range("E17").Value = '  -5.58%  '

$ws.Range("D18").Value = '90.95'
$ws.Range("E18").Value = '  -10.14%  '

$ws.Range("D19").Value = '0.06627'
$ws.Range("E19").Value = '  -1.50%  '

$ws.Range("D20").Value = '19.15'
$ws.Range("E20").Value = '  -9.36%  '

$ws.Range("D21").Value = '1.012'
$ws.Range("E21").Value = '  +1.05%  '

$ws.Range("D22").Value = '5.958'
$ws.Range("E22").Value = '  -6.72%  '

$ws.Range("D23").Value = '29.079.37'
$ws.Range("E23").Value = '  -3.97%  '

$ws.Range("D24").Value = '11.92'
$ws.Range("E24").Value = '  -4.11%  '

$ws.Range("D25").Value = '2.290'
$ws.Range("E25").Value = '  -1.09%  '

$ws.Range("B26").Value = 'EthereumClassic'
$ws.Range("C26").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D26").Value = '20.57'
$ws.Range("E26").Value = '  -5.74%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").Value = '156.00'
$ws.Range("E27").Value = '  -4.00%  '

$ws.Range("D28").Value = '6.208'
$ws.Range("E28").Value = '  -11.17%  '

$ws.Range("D29").Value = '2.256'
$ws.Range("E29").Value = '  -10.48%  '

$ws.Range("D30").Value = '126.68'
$ws.Range("E30").Value = '  -5.32%  '

$ws.Range("D31").Value = '1.040'
$ws.Range("E31").Value = '  -8.27%  '

$ws.Range("D32").Value = '0.09814'
$ws.Range("E32").Value = '  -6.82%  '

$ws.Range("D33").Value = '1.526'
$ws.Range("E33").Value = '  -9.60%  '

$ws.Range("D34").Value = '5.793'
$ws.Range("E34").Value = '  -7.38%  '

$ws.Range("E35").Value = '  -4.61%  '

$ws.Range("D36").Value = '0.02416'
$ws.Range("E36").Value = '  -7.81%  '

$ws.Range("D37").Value = '8.955'
$ws.Range("E37").Value = '  -11.13%  '

$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '1.286'
$ws.Range("E38").Value = '  -4.94%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.06301'
$ws.Range("E39").Value = '  -6.88%  '

$ws.Range("D40").Value = '0.6417'
$ws.Range("E40").Value = '  -7.87%  '

$ws.Range("D41").Value = '11.40'
$ws.Range("E41").Value = '  -10.04%  '

$ws.Range("D42").Value = '0.1985'
$ws.Range("E42").Value = '  -10.42%  '

$ws.Range("D43").Value = '1.009'
$ws.Range("E43").Value = '  +0.78%  '

$ws.Range("D44").Value = '0.6170'
$ws.Range("E44").Value = '  -8.98%  '

$ws.Range("D45").Value = '13.43'
$ws.Range("E45").Value = '  -5.93%  '

$ws.Range("D46").Value = '2.167'
$ws.Range("E46").Value = '  -7.08%  '

$ws.Range("D47").Value = '1.282'
$ws.Range("E47").Value = '  -2.24%  '

$ws.Range("D48").Value = '3.475'
$ws.Range("E48").Value = '  -4.50%  '

$ws.Range("D49").Value = '0.00000000327'
$ws.Range("E49").Value = '  -3.93%  '

$ws.Range("D50").Value = '0.06863'
$ws.Range("E50").Value = '  -5.98%  '

$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").Value = '1.115'
$ws.Range("E51").Value = '  -7.98%  '

# Restore default (Normal) cell style on the forced-text cells so no
# stray number-format styling is left behind on save.
foreach ($addr in $forceTextCells) {
    $ws.Range($addr).Style = "Normal"
}